$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 22 de Septiembre de 2020 a las 13:02"

# Swap Montserrat / Islas Malvinas rows (214 <-> 215 country names)
$ws.Cells.Item(214, 1).Value = "Islas Malvinas"
$ws.Cells.Item(215, 1).Value = "Montserrat"

# Update statistic columns (B=Casos totales, C=Nuevos casos, D=Casos activos,
# E=Recuperados, F=Casos criticos, G=Muertes hoy, H=Muertes)

# Row 4
$ws.Cells.Item(4, 2).Value = 7046444
$ws.Cells.Item(4, 3).Value = 228
$ws.Cells.Item(4, 4).Value = 4300731
$ws.Cells.Item(4, 5).Value = 2541198
$ws.Cells.Item(4, 7).Value = 9
$ws.Cells.Item(4, 8).Value = 204515

# Row 5
$ws.Cells.Item(5, 2).Value = 5567126
$ws.Cells.Item(5, 3).Value = 7021
$ws.Cells.Item(5, 5).Value = 980270
$ws.Cells.Item(5, 7).Value = 24
$ws.Cells.Item(5, 8).Value = 88989

# Row 16
$ws.Cells.Item(16, 2).Value = 429193
$ws.Cells.Item(16, 3).Value = 3712
$ws.Cells.Item(16, 4).Value = 363737
$ws.Cells.Item(16, 5).Value = 40800
$ws.Cells.Item(16, 7).Value = 178
$ws.Cells.Item(16, 8).Value = 24656

# Row 27
$ws.Cells.Item(27, 2).Value = 193374
$ws.Cells.Item(27, 3).Value = 2445
$ws.Cells.Item(27, 4).Value = 140751
$ws.Cells.Item(27, 5).Value = 51338
$ws.Cells.Item(27, 7).Value = 12
$ws.Cells.Item(27, 8).Value = 1285

# Row 33
$ws.Cells.Item(33, 2).Value = 114648
$ws.Cells.Item(33, 3).Value = 1059
$ws.Cells.Item(33, 4).Value = 92169
$ws.Cells.Item(33, 5).Value = 17976
$ws.Cells.Item(33, 7).Value = 45
$ws.Cells.Item(33, 8).Value = 4503

# Row 49
$ws.Cells.Item(49, 2).Value = 76104
$ws.Cells.Item(49, 3).Value = 206
$ws.Cells.Item(49, 4).Value = 73386
$ws.Cells.Item(49, 5).Value = 1927
$ws.Cells.Item(49, 7).Value = 6
$ws.Cells.Item(49, 8).Value = 791

# Row 69
$ws.Cells.Item(69, 2).Value = 39096
$ws.Cells.Item(69, 3).Value = 22
$ws.Cells.Item(69, 5).Value = 5075
$ws.Cells.Item(69, 7).Value = 1
$ws.Cells.Item(69, 8).Value = 1445

# Row 78
$ws.Cells.Item(78, 4).Value = 24218
$ws.Cells.Item(78, 5).Value = 1870

# Row 90
$ws.Cells.Item(90, 2).Value = 14759
$ws.Cells.Item(90, 3).Value = 21
$ws.Cells.Item(90, 4).Value = 11621
$ws.Cells.Item(90, 5).Value = 2836

# Row 97
$ws.Cells.Item(97, 2).Value = 10523
$ws.Cells.Item(97, 3).Value = 4
$ws.Cells.Item(97, 4).Value = 10011
$ws.Cells.Item(97, 5).Value = 241

# Row 145
$ws.Cells.Item(145, 2).Value = 2814
$ws.Cells.Item(145, 3).Value = 38
$ws.Cells.Item(145, 4).Value = 2113
$ws.Cells.Item(145, 5).Value = 678

# Row 179
$ws.Cells.Item(179, 2).Value = 448
$ws.Cells.Item(179, 3).Value = 11
$ws.Cells.Item(179, 4).Value = 414
$ws.Cells.Item(179, 5).Value = 34

# Row 182
$ws.Cells.Item(182, 2).Value = 355
$ws.Cells.Item(182, 3).Value = 5
$ws.Cells.Item(182, 5).Value = 32

# Row 214
$ws.Cells.Item(214, 4).Value = 13
$ws.Cells.Item(214, 8).Value = 0

# Row 215
$ws.Cells.Item(215, 4).Value = 12
$ws.Cells.Item(215, 8).Value = 1
